$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select B3 first (mirrors the authoring session's cursor position)
$ws.Range("B3").Select()

# Mark the task in row 3 ("End the task" column) as Done
$ws.Range("C3").Value = "Done"
